$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = 'Inflammatory-Mac'
$ws.Cells.Item(2,2).Value = 'Cd28'
$ws.Cells.Item(2,3).Value = 'Cd86'
$ws.Cells.Item(2,4).Value = 'ECs'
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 1.225147333333333
$ws.Cells.Item(2,8).Value = 3.675442
$ws.Cells.Item(2,9).Value = 0.2944933560673559
$ws.Cells.Item(2,10).Value = 0.2944933560673559
$ws.Cells.Item(2,11).Value = 1
$ws.Cells.Item(2,12).Value = 0.3333333333333333
$ws.Cells.Item(2,13).Value = 0.483405
$ws.Cells.Item(2,14).Value = 1.450215
$ws.Cells.Item(2,15).Value = 0.003094960828876145
$ws.Cells.Item(2,16).Value = 0.003094960828876144
$ws.Cells.Item(2,17).Value = 0.5922423466700001
$ws.Cells.Item(2,18).Value = 5.330181120030001
$ws.Cells.Item(2,19).Value = 0.0009114454013927415
$ws.Cells.Item(2,20).Value = 0.0009114454013927413

# Row 3
$ws.Cells.Item(3,1).Value = 'Inflammatory-Mac'
$ws.Cells.Item(3,2).Value = 'Cd28'
$ws.Cells.Item(3,3).Value = 'Cd86'
$ws.Cells.Item(3,4).Value = 'Inflammatory-Mac'
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 1.225147333333333
$ws.Cells.Item(3,8).Value = 3.675442
$ws.Cells.Item(3,9).Value = 0.2944933560673559
$ws.Cells.Item(3,10).Value = 0.2944933560673559
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 105.7018663333333
$ws.Cells.Item(3,14).Value = 317.105599
$ws.Cells.Item(3,15).Value = 0.676747521934545
$ws.Cells.Item(3,16).Value = 0.6767475219345449
$ws.Cells.Item(3,17).Value = 129.5003596666398
$ws.Cells.Item(3,18).Value = 1165.503236999758
$ws.Cells.Item(3,19).Value = 0.1992976489447707
$ws.Cells.Item(3,20).Value = 0.1992976489447707

# Row 4
$ws.Cells.Item(4,1).Value = 'Inflammatory-Mac'
$ws.Cells.Item(4,2).Value = 'Cd28'
$ws.Cells.Item(4,3).Value = 'Cd86'
$ws.Cells.Item(4,4).Value = 'Neutrophils'
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 1.225147333333333
$ws.Cells.Item(4,8).Value = 3.675442
$ws.Cells.Item(4,9).Value = 0.2944933560673559
$ws.Cells.Item(4,10).Value = 0.2944933560673559
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 18.83134933333333
$ws.Cells.Item(4,14).Value = 56.494048
$ws.Cells.Item(4,15).Value = 0.1205661682058513
$ws.Cells.Item(4,16).Value = 0.1205661682058513
$ws.Cells.Item(4,17).Value = 23.07117741880178
$ws.Cells.Item(4,18).Value = 207.640596769216
$ws.Cells.Item(4,19).Value = 0.0355059355031225
$ws.Cells.Item(4,20).Value = 0.0355059355031225

# Row 5
$ws.Cells.Item(5,1).Value = 'Inflammatory-Mac'
$ws.Cells.Item(5,2).Value = 'Cd28'
$ws.Cells.Item(5,3).Value = 'Cd86'
$ws.Cells.Item(5,4).Value = 'Resolving-Mac'
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 1.225147333333333
$ws.Cells.Item(5,8).Value = 3.675442
$ws.Cells.Item(5,9).Value = 0.2944933560673559
$ws.Cells.Item(5,10).Value = 0.2944933560673559
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 31.17437066666666
$ws.Cells.Item(5,14).Value = 93.523112
$ws.Cells.Item(5,15).Value = 0.1995913490307275
$ws.Cells.Item(5,16).Value = 0.1995913490307275
$ws.Cells.Item(5,17).Value = 38.19319709061156
$ws.Cells.Item(5,18).Value = 343.738773815504
$ws.Cells.Item(5,19).Value = 0.05877832621806994
$ws.Cells.Item(5,20).Value = 0.05877832621806994

# Row 6
$ws.Cells.Item(6,1).Value = 'MuSCs'
$ws.Cells.Item(6,2).Value = 'Cd28'
$ws.Cells.Item(6,3).Value = 'Cd86'
$ws.Cells.Item(6,4).Value = 'ECs'
$ws.Cells.Item(6,5).Value = 1
$ws.Cells.Item(6,6).Value = 0.3333333333333333
$ws.Cells.Item(6,7).Value = 0.005333666666666667
$ws.Cells.Item(6,8).Value = 0.016001
$ws.Cells.Item(6,9).Value = 0.001282073881300198
$ws.Cells.Item(6,10).Value = 0.001282073881300198
$ws.Cells.Item(6,11).Value = 1
$ws.Cells.Item(6,12).Value = 0.3333333333333333
$ws.Cells.Item(6,13).Value = 0.483405
$ws.Cells.Item(6,14).Value = 1.450215
$ws.Cells.Item(6,15).Value = 0.003094960828876145
$ws.Cells.Item(6,16).Value = 0.003094960828876144
$ws.Cells.Item(6,17).Value = 0.002578321135
$ws.Cells.Item(6,18).Value = 0.023204890215
$ws.Cells.Item(6,19).Value = [double]"3.967968442349317E-06"
$ws.Cells.Item(6,20).Value = [double]"3.967968442349316E-06"

# Row 7
$ws.Cells.Item(7,1).Value = 'MuSCs'
$ws.Cells.Item(7,2).Value = 'Cd28'
$ws.Cells.Item(7,3).Value = 'Cd86'
$ws.Cells.Item(7,4).Value = 'Inflammatory-Mac'
$ws.Cells.Item(7,5).Value = 1
$ws.Cells.Item(7,6).Value = 0.3333333333333333
$ws.Cells.Item(7,7).Value = 0.005333666666666667
$ws.Cells.Item(7,8).Value = 0.016001
$ws.Cells.Item(7,9).Value = 0.001282073881300198
$ws.Cells.Item(7,10).Value = 0.001282073881300198
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 105.7018663333333
$ws.Cells.Item(7,14).Value = 317.105599
$ws.Cells.Item(7,15).Value = 0.676747521934545
$ws.Cells.Item(7,16).Value = 0.6767475219345449
$ws.Cells.Item(7,17).Value = 0.5637785210665556
$ws.Cells.Item(7,18).Value = 5.074006689599
$ws.Cells.Item(7,19).Value = 0.000867640322106913
$ws.Cells.Item(7,20).Value = 0.0008676403221069129

# Row 8
$ws.Cells.Item(8,1).Value = 'MuSCs'
$ws.Cells.Item(8,2).Value = 'Cd28'
$ws.Cells.Item(8,3).Value = 'Cd86'
$ws.Cells.Item(8,4).Value = 'Neutrophils'
$ws.Cells.Item(8,5).Value = 1
$ws.Cells.Item(8,6).Value = 0.3333333333333333
$ws.Cells.Item(8,7).Value = 0.005333666666666667
$ws.Cells.Item(8,8).Value = 0.016001
$ws.Cells.Item(8,9).Value = 0.001282073881300198
$ws.Cells.Item(8,10).Value = 0.001282073881300198
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 18.83134933333333
$ws.Cells.Item(8,14).Value = 56.494048
$ws.Cells.Item(8,15).Value = 0.1205661682058513
$ws.Cells.Item(8,16).Value = 0.1205661682058513
$ws.Cells.Item(8,17).Value = 0.1004401402275556
$ws.Cells.Item(8,18).Value = 0.9039612620480001
$ws.Cells.Item(8,19).Value = 0.0001545747352251683
$ws.Cells.Item(8,20).Value = 0.0001545747352251683

# Row 9
$ws.Cells.Item(9,1).Value = 'MuSCs'
$ws.Cells.Item(9,2).Value = 'Cd28'
$ws.Cells.Item(9,3).Value = 'Cd86'
$ws.Cells.Item(9,4).Value = 'Resolving-Mac'
$ws.Cells.Item(9,5).Value = 1
$ws.Cells.Item(9,6).Value = 0.3333333333333333
$ws.Cells.Item(9,7).Value = 0.005333666666666667
$ws.Cells.Item(9,8).Value = 0.016001
$ws.Cells.Item(9,9).Value = 0.001282073881300198
$ws.Cells.Item(9,10).Value = 0.001282073881300198
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 31.17437066666666
$ws.Cells.Item(9,14).Value = 93.523112
$ws.Cells.Item(9,15).Value = 0.1995913490307275
$ws.Cells.Item(9,16).Value = 0.1995913490307275
$ws.Cells.Item(9,17).Value = 0.1662737016791111
$ws.Cells.Item(9,18).Value = 1.496463315112
$ws.Cells.Item(9,19).Value = 0.0002558908555257673
$ws.Cells.Item(9,20).Value = 0.0002558908555257673

# Row 10
$ws.Cells.Item(10,1).Value = 'Neutrophils'
$ws.Cells.Item(10,2).Value = 'Cd28'
$ws.Cells.Item(10,3).Value = 'Cd86'
$ws.Cells.Item(10,4).Value = 'ECs'
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 0.4766303333333334
$ws.Cells.Item(10,8).Value = 1.429891
$ws.Cells.Item(10,9).Value = 0.1145694584217375
$ws.Cells.Item(10,10).Value = 0.1145694584217375
$ws.Cells.Item(10,11).Value = 1
$ws.Cells.Item(10,12).Value = 0.3333333333333333
$ws.Cells.Item(10,13).Value = 0.483405
$ws.Cells.Item(10,14).Value = 1.450215
$ws.Cells.Item(10,15).Value = 0.003094960828876145
$ws.Cells.Item(10,16).Value = 0.003094960828876144
$ws.Cells.Item(10,17).Value = 0.2304054862850001
$ws.Cells.Item(10,18).Value = 2.073649376565001
$ws.Cells.Item(10,19).Value = 0.0003545879860008317
$ws.Cells.Item(10,20).Value = 0.0003545879860008316

# Row 11
$ws.Cells.Item(11,1).Value = 'Neutrophils'
$ws.Cells.Item(11,2).Value = 'Cd28'
$ws.Cells.Item(11,3).Value = 'Cd86'
$ws.Cells.Item(11,4).Value = 'Inflammatory-Mac'
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 0.4766303333333334
$ws.Cells.Item(11,8).Value = 1.429891
$ws.Cells.Item(11,9).Value = 0.1145694584217375
$ws.Cells.Item(11,10).Value = 0.1145694584217375
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 105.7018663333333
$ws.Cells.Item(11,14).Value = 317.105599
$ws.Cells.Item(11,15).Value = 0.676747521934545
$ws.Cells.Item(11,16).Value = 0.6767475219345449
$ws.Cells.Item(11,17).Value = 50.38071578441212
$ws.Cells.Item(11,18).Value = 453.4264420597091
$ws.Cells.Item(11,19).Value = 0.07753459707629373
$ws.Cells.Item(11,20).Value = 0.07753459707629373

# Row 12
$ws.Cells.Item(12,1).Value = 'Neutrophils'
$ws.Cells.Item(12,2).Value = 'Cd28'
$ws.Cells.Item(12,3).Value = 'Cd86'
$ws.Cells.Item(12,4).Value = 'Neutrophils'
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 0.4766303333333334
$ws.Cells.Item(12,8).Value = 1.429891
$ws.Cells.Item(12,9).Value = 0.1145694584217375
$ws.Cells.Item(12,10).Value = 0.1145694584217375
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 18.83134933333333
$ws.Cells.Item(12,14).Value = 56.494048
$ws.Cells.Item(12,15).Value = 0.1205661682058513
$ws.Cells.Item(12,16).Value = 0.1205661682058513
$ws.Cells.Item(12,17).Value = 8.975592309863112
$ws.Cells.Item(12,18).Value = 80.78033078876801
$ws.Cells.Item(12,19).Value = 0.01381320059532849
$ws.Cells.Item(12,20).Value = 0.01381320059532849

# Row 13
$ws.Cells.Item(13,1).Value = 'Neutrophils'
$ws.Cells.Item(13,2).Value = 'Cd28'
$ws.Cells.Item(13,3).Value = 'Cd86'
$ws.Cells.Item(13,4).Value = 'Resolving-Mac'
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 0.4766303333333334
$ws.Cells.Item(13,8).Value = 1.429891
$ws.Cells.Item(13,9).Value = 0.1145694584217375
$ws.Cells.Item(13,10).Value = 0.1145694584217375
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 31.17437066666666
$ws.Cells.Item(13,14).Value = 93.523112
$ws.Cells.Item(13,15).Value = 0.1995913490307275
$ws.Cells.Item(13,16).Value = 0.1995913490307275
$ws.Cells.Item(13,17).Value = 14.85865068231023
$ws.Cells.Item(13,18).Value = 133.727856140792
$ws.Cells.Item(13,19).Value = 0.02286707276411443
$ws.Cells.Item(13,20).Value = 0.02286707276411443

# Row 14
$ws.Cells.Item(14,1).Value = 'Resolving-Mac'
$ws.Cells.Item(14,2).Value = 'Cd28'
$ws.Cells.Item(14,3).Value = 'Cd86'
$ws.Cells.Item(14,4).Value = 'ECs'
$ws.Cells.Item(14,5).Value = 3
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(14,7).Value = 2.453075333333333
$ws.Cells.Item(14,8).Value = 7.359226
$ws.Cells.Item(14,9).Value = 0.5896551116296064
$ws.Cells.Item(14,10).Value = 0.5896551116296064
$ws.Cells.Item(14,11).Value = 1
$ws.Cells.Item(14,12).Value = 0.3333333333333333
$ws.Cells.Item(14,13).Value = 0.483405
$ws.Cells.Item(14,14).Value = 1.450215
$ws.Cells.Item(14,15).Value = 0.003094960828876145
$ws.Cells.Item(14,16).Value = 0.003094960828876144
$ws.Cells.Item(14,17).Value = 1.18582888151
$ws.Cells.Item(14,18).Value = 10.67245993359
$ws.Cells.Item(14,19).Value = 0.001824959473040222
$ws.Cells.Item(14,20).Value = 0.001824959473040222

# Row 15
$ws.Cells.Item(15,1).Value = 'Resolving-Mac'
$ws.Cells.Item(15,2).Value = 'Cd28'
$ws.Cells.Item(15,3).Value = 'Cd86'
$ws.Cells.Item(15,4).Value = 'Inflammatory-Mac'
$ws.Cells.Item(15,5).Value = 3
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(15,7).Value = 2.453075333333333
$ws.Cells.Item(15,8).Value = 7.359226
$ws.Cells.Item(15,9).Value = 0.5896551116296064
$ws.Cells.Item(15,10).Value = 0.5896551116296064
$ws.Cells.Item(15,11).Value = 3
$ws.Cells.Item(15,12).Value = 1
$ws.Cells.Item(15,13).Value = 105.7018663333333
$ws.Cells.Item(15,14).Value = 317.105599
$ws.Cells.Item(15,15).Value = 0.676747521934545
$ws.Cells.Item(15,16).Value = 0.6767475219345449
$ws.Cells.Item(15,17).Value = 259.2946409895971
$ws.Cells.Item(15,18).Value = 2333.651768906374
$ws.Cells.Item(15,19).Value = 0.3990476355913736
$ws.Cells.Item(15,20).Value = 0.3990476355913736

# Row 16
$ws.Cells.Item(16,1).Value = 'Resolving-Mac'
$ws.Cells.Item(16,2).Value = 'Cd28'
$ws.Cells.Item(16,3).Value = 'Cd86'
$ws.Cells.Item(16,4).Value = 'Neutrophils'
$ws.Cells.Item(16,5).Value = 3
$ws.Cells.Item(16,6).Value = 1
$ws.Cells.Item(16,7).Value = 2.453075333333333
$ws.Cells.Item(16,8).Value = 7.359226
$ws.Cells.Item(16,9).Value = 0.5896551116296064
$ws.Cells.Item(16,10).Value = 0.5896551116296064
$ws.Cells.Item(16,11).Value = 3
$ws.Cells.Item(16,12).Value = 1
$ws.Cells.Item(16,13).Value = 18.83134933333333
$ws.Cells.Item(16,14).Value = 56.494048
$ws.Cells.Item(16,15).Value = 0.1205661682058513
$ws.Cells.Item(16,16).Value = 0.1205661682058513
$ws.Cells.Item(16,17).Value = 46.19471854298311
$ws.Cells.Item(16,18).Value = 415.752466886848
$ws.Cells.Item(16,19).Value = 0.07109245737217516
$ws.Cells.Item(16,20).Value = 0.07109245737217515

# Row 17
$ws.Cells.Item(17,1).Value = 'Resolving-Mac'
$ws.Cells.Item(17,2).Value = 'Cd28'
$ws.Cells.Item(17,3).Value = 'Cd86'
$ws.Cells.Item(17,4).Value = 'Resolving-Mac'
$ws.Cells.Item(17,5).Value = 3
$ws.Cells.Item(17,6).Value = 1
$ws.Cells.Item(17,7).Value = 2.453075333333333
$ws.Cells.Item(17,8).Value = 7.359226
$ws.Cells.Item(17,9).Value = 0.5896551116296064
$ws.Cells.Item(17,10).Value = 0.5896551116296064
$ws.Cells.Item(17,11).Value = 3
$ws.Cells.Item(17,12).Value = 1
$ws.Cells.Item(17,13).Value = 31.17437066666666
$ws.Cells.Item(17,14).Value = 93.523112
$ws.Cells.Item(17,15).Value = 0.1995913490307275
$ws.Cells.Item(17,16).Value = 0.1995913490307275
$ws.Cells.Item(17,17).Value = 76.47307971459023
$ws.Cells.Item(17,18).Value = 688.257717431312
$ws.Cells.Item(17,19).Value = 0.1176900591930173
$ws.Cells.Item(17,20).Value = 0.1176900591930173
